$d = $word.ActiveDocument

# 1. Replace the last paragraph ("Adhere" Russian gloss) with the new block of
#    paragraphs: the same gloss (now carrying the "transcription" paragraph-mark
#    style) followed by the Treat / Threat / Premise word groups.
$target = $d.Paragraphs.Last.Range
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 wp14"><w:body><w:p w:rsidR="00963C35" w:rsidRPr="00963C35" w:rsidRDefault="00963C35" w:rsidP="00116A92"><w:pPr><w:spacing w:after="0"/><w:rPr><w:rStyle w:val="transcription"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rStyle w:val="transcription"/></w:rPr><w:t>Придерживаться, прилипать, твердо придерживаться</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/><w:rPr><w:rStyle w:val="transcription"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0"/><w:rPr><w:rStyle w:val="transcription"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rStyle w:val="transcription"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Treat </w:t></w:r><w:r><w:rPr><w:rStyle w:val="transcription"/></w:rPr><w:t>|triːt|</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/><w:rPr><w:rStyle w:val="transcription"/></w:rPr></w:pPr><w:r><w:rPr><w:rStyle w:val="transcription"/></w:rPr><w:t>Удовольствие, угощение, обращаться, отноститься, лечить, обрабатывать</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/><w:rPr><w:rStyle w:val="transcription"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0"/><w:rPr><w:rStyle w:val="transcription"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr><w:t>Threat</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rStyle w:val="transcription"/></w:rPr><w:t>|θret|</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/><w:rPr><w:rStyle w:val="transcription"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rStyle w:val="transcription"/></w:rPr><w:t>Угроза, опасность</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/><w:rPr><w:rStyle w:val="transcription"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0"/><w:rPr><w:rStyle w:val="transcription"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr><w:t>Premise</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rStyle w:val="transcription"/></w:rPr><w:t>|ˈpremɪs|</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rStyle w:val="transcription"/></w:rPr><w:t xml:space="preserve">Предпосылка, предпосылать, </w:t></w:r><w:r><w:t>недвижимость, дом с прилегающими постройками и участком</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>собственность, подлежащая передаче</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target.InsertXML($xml)

# InsertXML leaves a stray empty trailing paragraph behind (it carried the old
# paragraph-mark formatting of the paragraph that was replaced) - drop it.
$trailing = $d.Paragraphs.Last
if ($trailing.Range.Text -eq "`r") {
    $trailing.Range.Delete()
}

# InsertXML also drops the run-level "transcription" character style (rStyle),
# even though it keeps everything else - reapply it by locating each affected
# run through Find and setting its Style.
$styledRuns = @(
    "Придерживаться, прилипать, твердо придерживаться",
    "Treat ",
    "|triːt|",
    "Удовольствие, угощение, обращаться, отноститься, лечить, обрабатывать",
    "|θret|",
    "Угроза, опасность",
    "|ˈpremɪs|",
    "Предпосылка, предпосылать, "
)
foreach ($needle in $styledRuns) {
    $fr = $d.Content
    $ok = $fr.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($ok) {
        $fr.Style = "transcription"
    }
}
